$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-10-17"

# Update the "October (through 10-16)" label to "October (through 10-17)"
$ws.Range("A11").Value = "October (through 10-17)"

# Update the updated figures for row 10 (September), row 11 (October), row 12 (Total)
$ws.Range("H10").Value = 178

$ws.Range("B11").Value = 14
$ws.Range("E11").Value = 45
$ws.Range("G11").Value = 82
$ws.Range("H11").Value = 105

$ws.Range("B12").Value = 240
$ws.Range("E12").Value = 593
$ws.Range("G12").Value = 983
$ws.Range("H12").Value = 1354
